$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2916.1018
$ws.Range("I138").Value = 1948.1482
$ws.Range("J138").Value = 3732.8125
$ws.Range("K138").Value = 5844.444600000001
$ws.Range("L138").Value = 11198.4375
$ws.Range("M138").Value = -704.4446000000007
$ws.Range("N138").Value = -21478.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9338.846
$ws.Range("I2").Value = 783.6667
$ws.Range("K2").Value = 783.6667
$ws.Range("M2").Value = -670.6667

$ws.Range("H32").Value = 3658.7297
$ws.Range("I32").Value = 2936.7646
$ws.Range("J32").Value = 11841
$ws.Range("K32").Value = 2936.7646
$ws.Range("L32").Value = 11841
$ws.Range("M32").Value = -2649.7646
$ws.Range("N32").Value = -12415

$ws.Range("H45").Value = 76927660
$ws.Range("I45").Value = 125002010
$ws.Range("K45").Value = 125002010
$ws.Range("M45").Value = -125001633

$ws.Range("H61").Value = 5705.7856
$ws.Range("I61").Value = 5705.7856
$ws.Range("K61").Value = 5705.7856
$ws.Range("M61").Value = -5493.7856

$ws.Range("H74").Value = 15875118
$ws.Range("I74").Value = 16668400
$ws.Range("K74").Value = 16668400
$ws.Range("M74").Value = -16667526

$ws.Range("H77").Value = 15875118
$ws.Range("I77").Value = 16668400
$ws.Range("K77").Value = 83342000
$ws.Range("M77").Value = -83337632

$ws.Range("H102").Value = 1221.4166
$ws.Range("I102").Value = 1221.4166
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1221.4166
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 400.5834
$ws.Range("N102").ClearContents()

$ws.Range("H116").Value = 9338.846
$ws.Range("I116").Value = 783.6667
$ws.Range("K116").Value = 783.6667
$ws.Range("M116").Value = 1510.3333

$ws.Range("H122").Value = 3158.739
$ws.Range("I122").Value = 2093.7144
$ws.Range("K122").Value = 6281.1432
$ws.Range("M122").Value = -3831.1432

$ws.Range("H132").Value = 4927.636
$ws.Range("I132").Value = 2336.15
$ws.Range("J132").Value = 8914.538
$ws.Range("K132").Value = 7008.450000000001
$ws.Range("L132").Value = 26743.614
$ws.Range("M132").Value = -4478.450000000001
$ws.Range("N132").Value = -31803.614

$ws.Range("H136").Value = 5705.7856
$ws.Range("I136").Value = 5705.7856
$ws.Range("K136").Value = 17117.3568
$ws.Range("M136").Value = -14567.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9338.846
$ws.Range("I3").Value = 783.6667
$ws.Range("K3").Value = 783.6667
$ws.Range("M3").Value = -669.6667

$ws.Range("H20").Value = 5551.1113
$ws.Range("I20").Value = 5330.4
$ws.Range("K20").Value = 5330.4
$ws.Range("M20").Value = -5083.4

$ws.Range("H64").Value = 2291.125
$ws.Range("I64").Value = 1680.75
$ws.Range("J64").Value = 2901.5
$ws.Range("K64").Value = 1680.75
$ws.Range("L64").Value = 2901.5
$ws.Range("M64").Value = -1455.75
$ws.Range("N64").Value = -3351.5

$ws.Range("H67").Value = 2291.125
$ws.Range("I67").Value = 1680.75
$ws.Range("J67").Value = 2901.5
$ws.Range("K67").Value = 1680.75
$ws.Range("L67").Value = 2901.5
$ws.Range("M67").Value = -900.75
$ws.Range("N67").Value = -4461.5

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws.Range("H134").Value = 2268.0386
$ws.Range("I134").Value = 1415.6666
$ws.Range("K134").Value = 4246.9998
$ws.Range("M134").Value = -1711.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23842.2
$ws.Range("I31").Value = 3381.5264
$ws.Range("J31").Value = 69577.82000000001
$ws.Range("K31").Value = 3381.5264
$ws.Range("L31").Value = 69577.82000000001
$ws.Range("M31").Value = -3086.5264
$ws.Range("N31").Value = -70167.82000000001

$ws.Range("H34").Value = 23842.2
$ws.Range("I34").Value = 3381.5264
$ws.Range("J34").Value = 69577.82000000001
$ws.Range("K34").Value = 3381.5264
$ws.Range("L34").Value = 69577.82000000001
$ws.Range("M34").Value = -3179.5264
$ws.Range("N34").Value = -69981.82000000001

$ws.Range("H58").Value = 4267.1274
$ws.Range("I58").Value = 3435.578
$ws.Range("J58").Value = 8009.1
$ws.Range("K58").Value = 3435.578
$ws.Range("L58").Value = 8009.1
$ws.Range("M58").Value = -3232.578
$ws.Range("N58").Value = -8415.1

$ws.Range("H86").Value = 10535.777
$ws.Range("I86").Value = 5457
$ws.Range("J86").Value = 11986.857
$ws.Range("K86").Value = 5457
$ws.Range("L86").Value = 11986.857
$ws.Range("M86").Value = -4334
$ws.Range("N86").Value = -14232.857

$ws.Range("H89").Value = 10535.777
$ws.Range("I89").Value = 5457
$ws.Range("J89").Value = 11986.857
$ws.Range("K89").Value = 27285
$ws.Range("L89").Value = 59934.285
$ws.Range("M89").Value = -21669
$ws.Range("N89").Value = -71166.285

$ws.Range("H94").Value = 4574.6665
$ws.Range("J94").Value = 4089.6
$ws.Range("L94").Value = 4089.6
$ws.Range("N94").Value = -4991.6

$ws.Range("H107").Value = 1217.6364
$ws.Range("I107").Value = 1187.1538
$ws.Range("K107").Value = 1187.1538
$ws.Range("M107").Value = 732.8462

$ws.Range("H136").Value = 4267.1274
$ws.Range("I136").Value = 3435.578
$ws.Range("J136").Value = 8009.1
$ws.Range("K136").Value = 10306.734
$ws.Range("L136").Value = 24027.3
$ws.Range("M136").Value = -7756.734
$ws.Range("N136").Value = -29127.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 19026.5
$ws.Range("I95").Value = 19026
$ws.Range("K95").Value = 57078
$ws.Range("M95").Value = -55019

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H102").Value = 7999.8
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 7999.8
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 23999.4
$ws.Range("N102").Value = -28867.4
$ws.Range("M102").ClearContents()

$ws.Range("H125").Value = 14218.889
$ws.Range("I125").Value = 7666.6665
$ws.Range("J125").Value = 17495
$ws.Range("K125").Value = 22999.9995
$ws.Range("L125").Value = 52485
$ws.Range("M125").Value = -18079.9995
$ws.Range("N125").Value = -62325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15668.685
$ws.Range("I70").Value = 5227.3335
$ws.Range("J70").Value = 33568.145
$ws.Range("K70").Value = 5227.3335
$ws.Range("L70").Value = 33568.145
$ws.Range("M70").Value = -4957.3335
$ws.Range("N70").Value = -34108.145

$ws.Range("H73").Value = 15668.685
$ws.Range("I73").Value = 5227.3335
$ws.Range("J73").Value = 33568.145
$ws.Range("K73").Value = 5227.3335
$ws.Range("L73").Value = 33568.145
$ws.Range("M73").Value = -4291.3335
$ws.Range("N73").Value = -35440.145

$ws.Range("H102").Value = 2791.7
$ws.Range("I102").Value = 957.8570999999999
$ws.Range("J102").Value = 7070.6665
$ws.Range("K102").Value = 957.8570999999999
$ws.Range("L102").Value = 7070.6665
$ws.Range("M102").Value = 664.1429000000001
$ws.Range("N102").Value = -10314.6665

$ws.Range("H126").Value = 4645.2856
$ws.Range("J126").Value = 6901.625
$ws.Range("L126").Value = 20704.875
$ws.Range("N126").Value = -25644.875

$ws.Range("H132").Value = 58023.6
$ws.Range("I132").Value = 62581
$ws.Range("K132").Value = 187743
$ws.Range("M132").Value = -185213

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 21424
$ws.Range("J43").Value = 27987.4
$ws.Range("L43").Value = 27987.4
$ws.Range("N43").Value = -28373.4

$ws.Range("H55").Value = 2174505.8
$ws.Range("I55").Value = 3125535
$ws.Range("J55").Value = 725
$ws.Range("K55").Value = 3125535
$ws.Range("L55").Value = 725
$ws.Range("M55").Value = -3125362
$ws.Range("N55").Value = -1071

$ws.Range("H100").Value = 4500.0835
$ws.Range("I100").Value = 2399.7
$ws.Range("J100").Value = 15002
$ws.Range("K100").Value = 2399.7
$ws.Range("L100").Value = 15002
$ws.Range("M100").Value = -1858.7
$ws.Range("N100").Value = -16084

$ws.Range("H132").Value = 3609.5
$ws.Range("I132").Value = 2312.2964
$ws.Range("J132").Value = 7501.1113
$ws.Range("K132").Value = 6936.889200000001
$ws.Range("L132").Value = 22503.3339
$ws.Range("M132").Value = -4406.889200000001
$ws.Range("N132").Value = -27563.3339

$ws.Range("H136").Value = 4797.263
$ws.Range("I136").Value = 2210.6155
$ws.Range("J136").Value = 10401.667
$ws.Range("K136").Value = 6631.8465
$ws.Range("L136").Value = 31205.001
$ws.Range("M136").Value = -4081.8465
$ws.Range("N136").Value = -36305.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3693.5789
$ws.Range("I136").Value = 2456.58
$ws.Range("J136").Value = 12529.286
$ws.Range("K136").Value = 7369.74
$ws.Range("L136").Value = 37587.858
$ws.Range("M136").Value = -4819.74
$ws.Range("N136").Value = -42687.858
